# "updated activity till excel form"
# Rewrites the per-match batting stats (runs/balls/fours/sixes) for the
# existing 13 rows and appends a new 14th match as row 15.
# Values are entered with a leading apostrophe so Excel stores them as
# literal text (matching the sheet's existing text-as-number convention)
# instead of auto-converting the numeric-looking strings to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 ---
$ws.Cells.Item(2,3).Value = "'31"
$ws.Cells.Item(2,4).Value = "'20"
$ws.Cells.Item(2,5).Value = "'5"
$ws.Cells.Item(2,6).Value = "'0"

# --- row 3 ---
$ws.Cells.Item(3,3).Value = "'19"
$ws.Cells.Item(3,4).Value = "'15"
$ws.Cells.Item(3,5).Value = "'2"
$ws.Cells.Item(3,6).Value = "'0"

# --- row 4 ---
$ws.Cells.Item(4,3).Value = "'11"
$ws.Cells.Item(4,4).Value = "'8"
$ws.Cells.Item(4,5).Value = "'1"
$ws.Cells.Item(4,6).Value = "'1"

# --- row 5 ---
$ws.Cells.Item(5,3).Value = "'6"
$ws.Cells.Item(5,4).Value = "'7"
$ws.Cells.Item(5,5).Value = "'1"
$ws.Cells.Item(5,6).Value = "'0"

# --- row 6 ---
$ws.Cells.Item(6,3).Value = "'69"
$ws.Cells.Item(6,4).Value = "'47"
$ws.Cells.Item(6,5).Value = "'4"
$ws.Cells.Item(6,6).Value = "'4"

# --- row 7 ---
$ws.Cells.Item(7,3).Value = "'4"
$ws.Cells.Item(7,4).Value = "'4"
$ws.Cells.Item(7,5).Value = "'1"
$ws.Cells.Item(7,6).Value = "'0"

# --- row 8 ---
$ws.Cells.Item(8,3).Value = "'26"
$ws.Cells.Item(8,4).Value = "'34"
$ws.Cells.Item(8,5).Value = "'2"
$ws.Cells.Item(8,6).Value = "'0"

# --- row 9 ---
$ws.Cells.Item(9,3).Value = "'1"
$ws.Cells.Item(9,4).Value = "'4"
$ws.Cells.Item(9,5).Value = "'0"
$ws.Cells.Item(9,6).Value = "'0"

# --- row 10 ---
$ws.Cells.Item(10,3).Value = "'5"
$ws.Cells.Item(10,4).Value = "'6"
$ws.Cells.Item(10,5).Value = "'0"
$ws.Cells.Item(10,6).Value = "'0"

# --- row 11 ---
$ws.Cells.Item(11,3).Value = "'5"
$ws.Cells.Item(11,4).Value = "'5"
$ws.Cells.Item(11,5).Value = "'1"
$ws.Cells.Item(11,6).Value = "'0"

# --- row 12 ---
$ws.Cells.Item(12,3).Value = "'57"
$ws.Cells.Item(12,4).Value = "'36"
$ws.Cells.Item(12,5).Value = "'6"
$ws.Cells.Item(12,6).Value = "'1"

# --- row 13 ---
$ws.Cells.Item(13,3).Value = "'50"
$ws.Cells.Item(13,4).Value = "'27"
$ws.Cells.Item(13,5).Value = "'7"
$ws.Cells.Item(13,6).Value = "'2"

# --- row 14 ---
$ws.Cells.Item(14,3).Value = "'3"
$ws.Cells.Item(14,4).Value = "'7"
$ws.Cells.Item(14,5).Value = "'0"
$ws.Cells.Item(14,6).Value = "'0"

# --- row 15 (new match, appended) ---
$ws.Cells.Item(15,1).Value = "Steven Smith (c)"
$ws.Cells.Item(15,2).Value = "Rajasthan Royals"
$ws.Cells.Item(15,3).Value = "'24"
$ws.Cells.Item(15,4).Value = "'17"
$ws.Cells.Item(15,5).Value = "'2"
$ws.Cells.Item(15,6).Value = "'1"

# Keep the "numbers stored as text" warning suppressed over the full used
# range now that it has grown by one row (A1:F14 -> A1:F15).
$ws.Range("A1:F15").Errors.Item(9).Ignore = $true
